# Apply workbook edits per the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet/tab
$ws.Name = "overall-expense"

# Update row 2
$ws.Range("B2").Value = "pav bhaji"
$ws.Range("C2").Value = 100
$ws.Range("D2").Value = "Shubham More"
$ws.Range("E2").Value = "EQUAL"
$ws.Range("F2").Value = "Shubham More: 50; Dev: 50"

# Update row 3
$ws.Range("B3").Value = "pav bhaji"
$ws.Range("C3").Value = 200
$ws.Range("D3").Value = "Shubham More"
$ws.Range("E3").Value = "PERCENTAGE"
$ws.Range("F3").Value = "Shubham More: 60; Dev: 140"

# Update row 4
$ws.Range("B4").Value = "pav bhaji"
$ws.Range("C4").Value = 100
$ws.Range("D4").Value = "Shubham More"
$ws.Range("E4").Value = "EXACT"
$ws.Range("F4").Value = "Shubham More: 45; Dev: 55"
